$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")

# --- Populate new rows 15-29 (IAM array variable mapping for Physical Trough model) ---
# Cell values are written in the exact order the original author entered them,
# so that newly-introduced shared strings land in the same order as the target file.
$ws.Range('C15').Value = 'csp_dtr_sca_iam0_1'
$ws.Range('C16').Value = 'csp_dtr_sca_iam1_1'
$ws.Range('C18').Value = 'csp_dtr_sca_iam0_2'
$ws.Range('C17').Value = 'csp_dtr_sca_iam2_1'
$ws.Range('D15').Value = 'IAMs_1[0]'
$ws.Range('D16').Value = 'IAMs_1[1]'
$ws.Range('D17').Value = 'IAMs_1[2]'
$ws.Range('E15').Value = 'Physical Trough Collector Type 1'
$ws.Range('E18').Value = 'Physical Trough Collector Type 2'
$ws.Range('E21').Value = 'Physical Trough Collector Type 3'
$ws.Range('F15').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G15').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('C27').Value = 'IamF0'
$ws.Range('E27').Value = 'Physical Trough Collector Header'
$ws.Range('D27').Value = 'IAM_matrix'
$ws.Range('F27').Value = 'combining collector IAM coef. Arrays into 1 output matrix'
$ws.Range('C28').Value = 'IamF1'
$ws.Range('C29').Value = 'IamF2'
$ws.Range('C19').Value = 'csp_dtr_sca_iam1_2'
$ws.Range('C20').Value = 'csp_dtr_sca_iam2_2'
$ws.Range('D18').Value = 'IAMs_2[0]'
$ws.Range('D19').Value = 'IAMs_2[1]'
$ws.Range('D20').Value = 'IAMs_2[2]'
$ws.Range('C21').Value = 'csp_dtr_sca_iam0_3'
$ws.Range('C22').Value = 'csp_dtr_sca_iam1_3'
$ws.Range('C23').Value = 'csp_dtr_sca_iam2_3'
$ws.Range('D21').Value = 'IAMs_3[0]'
$ws.Range('D22').Value = 'IAMs_3[1]'
$ws.Range('D23').Value = 'IAMs_3[2]'
$ws.Range('C24').Value = 'csp_dtr_sca_iam0_4'
$ws.Range('C25').Value = 'csp_dtr_sca_iam1_4'
$ws.Range('C26').Value = 'csp_dtr_sca_iam2_4'
$ws.Range('D24').Value = 'IAMs_4[0]'
$ws.Range('D25').Value = 'IAMs_4[1]'
$ws.Range('D26').Value = 'IAMs_4[2]'
$ws.Range('E24').Value = 'Physical Trough Collector Type 4'
$ws.Range('A15').Value = 'Changed name'
$ws.Range('B15').Value = 'number'
$ws.Range('A16').Value = 'Changed name'
$ws.Range('B16').Value = 'number'
$ws.Range('E16').Value = 'Physical Trough Collector Type 1'
$ws.Range('F16').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G16').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A17').Value = 'Changed name'
$ws.Range('B17').Value = 'number'
$ws.Range('E17').Value = 'Physical Trough Collector Type 1'
$ws.Range('F17').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G17').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A18').Value = 'Changed name'
$ws.Range('B18').Value = 'number'
$ws.Range('F18').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G18').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A19').Value = 'Changed name'
$ws.Range('B19').Value = 'number'
$ws.Range('E19').Value = 'Physical Trough Collector Type 2'
$ws.Range('F19').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G19').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A20').Value = 'Changed name'
$ws.Range('B20').Value = 'number'
$ws.Range('E20').Value = 'Physical Trough Collector Type 2'
$ws.Range('F20').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G20').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A21').Value = 'Changed name'
$ws.Range('B21').Value = 'number'
$ws.Range('F21').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G21').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A22').Value = 'Changed name'
$ws.Range('B22').Value = 'number'
$ws.Range('E22').Value = 'Physical Trough Collector Type 3'
$ws.Range('F22').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G22').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A23').Value = 'Changed name'
$ws.Range('B23').Value = 'number'
$ws.Range('E23').Value = 'Physical Trough Collector Type 3'
$ws.Range('F23').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G23').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A24').Value = 'Changed name'
$ws.Range('B24').Value = 'number'
$ws.Range('F24').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G24').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A25').Value = 'Changed name'
$ws.Range('B25').Value = 'number'
$ws.Range('E25').Value = 'Physical Trough Collector Type 4'
$ws.Range('F25').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G25').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A26').Value = 'Changed name'
$ws.Range('B26').Value = 'number'
$ws.Range('E26').Value = 'Physical Trough Collector Type 4'
$ws.Range('F26').Value = 'Allowing table/array of IAM coefficients as inputs'
$ws.Range('G26').Value = 'Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables'
$ws.Range('A27').Value = 'Changed name'
$ws.Range('B27').Value = 'array'
$ws.Range('G27').Value = 'Ty'
$ws.Range('A28').Value = 'Changed name'
$ws.Range('B28').Value = 'array'
$ws.Range('D28').Value = 'IAM_matrix'
$ws.Range('E28').Value = 'Physical Trough Collector Header'
$ws.Range('F28').Value = 'combining collector IAM coef. Arrays into 1 output matrix'
$ws.Range('G28').Value = 'Ty'
$ws.Range('A29').Value = 'Changed name'
$ws.Range('B29').Value = 'array'
$ws.Range('D29').Value = 'IAM_matrix'
$ws.Range('E29').Value = 'Physical Trough Collector Header'
$ws.Range('F29').Value = 'combining collector IAM coef. Arrays into 1 output matrix'
$ws.Range('G29').Value = 'Ty'

# --- Column F is now much longer text, widen it (and drop the old bestFit flag) ---
# (50.14 chars round-trips through the engine's pixel quantization to exactly width=51)
$ws.Columns.Item(6).ColumnWidth = 50.14

# --- Extend the "Type" list validation down through the newly added rows ---
$rngValidation = $ws.Range("A2:A58")
$rngValidation.Validation.Delete()
$rngValidation.Validation.Add(3, 1, 1, "Types")
$rngValidation.Validation.IgnoreBlank = $false

# --- Move/resize the instructional callout box out of the way of the new rows ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 11953875 / 12700
$shp.Top = 4010025 / 12700
$shp.Width = 3990975 / 12700
$shp.Height = 1581150 / 12700

# --- Update the active selection to reflect where the author ended up editing ---
$ws.Range("D29").Select()
